$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted as row 70 (Berenjena / Femacal de La
# Calera), pushing the previously-existing rows 70..147 down to 71..148.
$ws.Rows.Item(70).Insert()

$ws.Cells.Item(70, 1).Value = 3
$ws.Cells.Item(70, 2).Value = 'Femacal de La Calera'
$ws.Cells.Item(70, 3).Value = 'Coquimbo'
$ws.Cells.Item(70, 4).Value = 44467
$ws.Cells.Item(70, 5).Value = 5
$ws.Cells.Item(70, 6).Value = 100112001
$ws.Cells.Item(70, 7).Value = 'Berenjena'
$ws.Cells.Item(70, 8).Value = 'Sin especificar'
$ws.Cells.Item(70, 9).Value = 'Primera'
$ws.Cells.Item(70, 10).Value = 110
$ws.Cells.Item(70, 11).Value = 9000
$ws.Cells.Item(70, 12).Value = 9500
$ws.Cells.Item(70, 13).Value = 9273
$ws.Cells.Item(70, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(70, 15).Value = 'Región de Arica y Parinacota'
$ws.Cells.Item(70, 16).Value = 155
$ws.Cells.Item(70, 17).Value = 60
$ws.Cells.Item(70, 18).Value = 'Hortaliza'
